$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C7:C14 numeric values (content column) per diff
$ws.Cells.Item(7, 3).Value = 5
$ws.Cells.Item(8, 3).Value = 8
$ws.Cells.Item(9, 3).Value = 9
$ws.Cells.Item(10, 3).Value = 11
$ws.Cells.Item(11, 3).Value = 13
$ws.Cells.Item(12, 3).Value = 16
$ws.Cells.Item(13, 3).Value = 17
$ws.Cells.Item(14, 3).Value = 18

# Row 15: btc -> ftm ; Significant Related Persons -> Associated Technology ; 12 -> "blockchain"
$ws.Cells.Item(15, 1).Value = "ftm"
$ws.Cells.Item(15, 2).Value = "Associated Technology"
$ws.Cells.Item(15, 3).Value = "blockchain"

# Row 16: Significant Related Persons -> Significant Persons ; 13 -> "elon musk"
$ws.Cells.Item(16, 2).Value = "Significant Persons"
$ws.Cells.Item(16, 3).Value = "elon musk"

# Row 17: btc -> eth ; 14 -> "Vitalik Buterin"
$ws.Cells.Item(17, 1).Value = "eth"
$ws.Cells.Item(17, 3).Value = "Vitalik Buterin"

# Row 18: 15 -> "vip1"
$ws.Cells.Item(18, 3).Value = "vip1"

# Row 19: 16 -> "vip2"
$ws.Cells.Item(19, 3).Value = "vip2"

# Row 20: 17 -> "vip3"
$ws.Cells.Item(20, 3).Value = "vip3"

# Row 21: Significant Related Persons -> Associated Technology ; 18 -> "tech1"
$ws.Cells.Item(21, 2).Value = "Associated Technology"
$ws.Cells.Item(21, 3).Value = "tech1"

# Row 22: ftm -> btc ; blockchain -> "tech2"
$ws.Cells.Item(22, 1).Value = "btc"
$ws.Cells.Item(22, 3).Value = "tech2"

# Row 23: Luna -> btc ; Correlated Concept -> Associated Technology ; btc -> "tech3"
$ws.Cells.Item(23, 1).Value = "btc"
$ws.Cells.Item(23, 2).Value = "Associated Technology"
$ws.Cells.Item(23, 3).Value = "tech3"

# Row 24: Significant Persons -> Correlated Concept ; elon musk -> "concept1"
$ws.Cells.Item(24, 2).Value = "Correlated Concept"
$ws.Cells.Item(24, 3).Value = "concept1"

# Row 25 (eth, Significant Related Persons, Vitalik Buterin, test111) is removed entirely
$ws.Rows.Item(25).Delete()
